$wb = $excel.ActiveWorkbook

# zh-cn sheet: update row 5 Handoff/Handback datetime cells
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-23 09:21:01"
$wsZhCn.Range("G5").Value = "2016-02-23 09:21:49"

# de-de sheet: update row 5 Handoff/Handback datetime cells
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-23 09:21:13"
$wsDeDe.Range("G5").Value = "2016-02-23 09:22:11"
